# B6-PowerPoint.pptx edit
# The three data tables (slides 14, 15 and 16) had their table style
# switched from the custom "Table_0" style to the built-in PowerPoint
# table style {5016DCA5-F7D0-478C-94C0-965BF381F7D7}.
#
# Walk every slide and every shape on it; wherever a shape carries a
# table, re-apply the new style id. This is robust to shape ordering
# and to the table not always being shape #1 on the slide.

$p = $ppt.ActivePresentation
$newStyleId = "{5016DCA5-F7D0-478C-94C0-965BF381F7D7}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -ne $newStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
